$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 570.8
$ws.Range("I103").Value = 602
$ws.Range("J103").Value = 550
$ws.Range("K103").Value = 1806
$ws.Range("L103").Value = 1650
$ws.Range("M103").Value = -1220
$ws.Range("N103").Value = -2822

$ws.Range("H111").Value = 2515.875
$ws.Range("I111").Value = 2565.5
$ws.Range("J111").Value = 2433.1667
$ws.Range("K111").Value = 7696.5
$ws.Range("L111").Value = 7299.500100000001
$ws.Range("M111").Value = -4629.5
$ws.Range("N111").Value = -13433.5001

$ws.Range("H132").Value = 4078.8076
$ws.Range("I132").Value = 3947.9333
$ws.Range("J132").Value = 4257.273
$ws.Range("K132").Value = 11843.7999
$ws.Range("L132").Value = 12771.819
$ws.Range("M132").Value = -9313.7999
$ws.Range("N132").Value = -17831.819

$ws.Range("H137").Value = 1193.7742
$ws.Range("I137").Value = 985.5417
$ws.Range("K137").Value = 2956.6251
$ws.Range("M137").Value = -406.6251000000002

$ws.Range("H138").Value = 2953.1077
$ws.Range("I138").Value = 1971
$ws.Range("J138").Value = 4872.6816
$ws.Range("K138").Value = 5913
$ws.Range("L138").Value = 14618.0448
$ws.Range("M138").Value = -773
$ws.Range("N138").Value = -24898.0448

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1826.4445
$ws.Range("I45").Value = 1544.8
$ws.Range("J45").Value = 2178.5
$ws.Range("K45").Value = 1544.8
$ws.Range("L45").Value = 2178.5
$ws.Range("M45").Value = -1167.8
$ws.Range("N45").Value = -2932.5

$ws.Range("H61").Value = 1059.3959
$ws.Range("I61").Value = 1035.225
$ws.Range("J61").Value = 1180.25
$ws.Range("K61").Value = 1035.225
$ws.Range("L61").Value = 1180.25
$ws.Range("M61").Value = -823.2249999999999
$ws.Range("N61").Value = -1604.25

$ws.Range("H74").Value = 1091.6757
$ws.Range("I74").Value = 1117.091
$ws.Range("J74").Value = 882
$ws.Range("K74").Value = 1117.091
$ws.Range("L74").Value = 882
$ws.Range("M74").Value = -243.0909999999999
$ws.Range("N74").Value = -2630

$ws.Range("H77").Value = 1091.6757
$ws.Range("I77").Value = 1117.091
$ws.Range("J77").Value = 882
$ws.Range("K77").Value = 5585.455
$ws.Range("L77").Value = 4410
$ws.Range("M77").Value = -1217.455
$ws.Range("N77").Value = -13146

$ws.Range("H131").Value = 32810
$ws.Range("J131").Value = 32810
$ws.Range("L131").Value = 32810
$ws.Range("N131").Value = -42890

$ws.Range("H136").Value = 1059.3959
$ws.Range("I136").Value = 1035.225
$ws.Range("J136").Value = 1180.25
$ws.Range("K136").Value = 3105.675
$ws.Range("L136").Value = 3540.75
$ws.Range("M136").Value = -555.6749999999997
$ws.Range("N136").Value = -8640.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 20000
$ws.Range("J19").Value = 20000
$ws.Range("L19").Value = 20000
$ws.Range("N19").Value = -20346

$ws.Range("H109").Value = 21403.092
$ws.Range("J109").Value = 21403.092
$ws.Range("L109").Value = 21403.092
$ws.Range("N109").Value = -24177.092

$ws.Range("H127").Value = 20140
$ws.Range("J127").Value = 20140
$ws.Range("L127").Value = 20140
$ws.Range("N127").Value = -30060

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 8821.412
$ws.Range("I4").Value = 6000
$ws.Range("J4").Value = 9997
$ws.Range("K4").Value = 6000
$ws.Range("L4").Value = 9997
$ws.Range("M4").Value = -5888
$ws.Range("N4").Value = -10221

$ws.Range("H22").Value = 5288.8
$ws.Range("I22").Value = 7364.5713
$ws.Range("J22").Value = 445.33334
$ws.Range("K22").Value = 7364.5713
$ws.Range("L22").Value = 445.33334
$ws.Range("M22").Value = -7014.5713
$ws.Range("N22").Value = -1145.33334

$ws.Range("H31").Value = 2101.5881
$ws.Range("I31").Value = 1439.4584
$ws.Range("J31").Value = 3690.7
$ws.Range("K31").Value = 1439.4584
$ws.Range("L31").Value = 3690.7
$ws.Range("M31").Value = -1144.4584
$ws.Range("N31").Value = -4280.7

$ws.Range("H34").Value = 2101.5881
$ws.Range("I34").Value = 1439.4584
$ws.Range("J34").Value = 3690.7
$ws.Range("K34").Value = 1439.4584
$ws.Range("L34").Value = 3690.7
$ws.Range("M34").Value = -1237.4584
$ws.Range("N34").Value = -4094.7

$ws.Range("H58").Value = 772781.5600000001
$ws.Range("I58").Value = 1323638
$ws.Range("J58").Value = 1582.45
$ws.Range("K58").Value = 1323638
$ws.Range("L58").Value = 1582.45
$ws.Range("M58").Value = -1323435
$ws.Range("N58").Value = -1988.45

$ws.Range("H132").Value = 484494.7
$ws.Range("I132").Value = 615819.9
$ws.Range("K132").Value = 1847459.7
$ws.Range("M132").Value = -1844929.7

$ws.Range("H134").Value = 1751.4615
$ws.Range("I134").Value = 1514.95
$ws.Range("J134").Value = 2539.8333
$ws.Range("K134").Value = 4544.85
$ws.Range("L134").Value = 7619.499899999999
$ws.Range("M134").Value = -2009.85
$ws.Range("N134").Value = -12689.4999

$ws.Range("H136").Value = 772781.5600000001
$ws.Range("I136").Value = 1323638
$ws.Range("J136").Value = 1582.45
$ws.Range("K136").Value = 3970914
$ws.Range("L136").Value = 4747.35
$ws.Range("M136").Value = -3968364
$ws.Range("N136").Value = -9847.35

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 76.583336
$ws.Range("I6").Value = 65.36364
$ws.Range("J6").Value = 200
$ws.Range("K6").Value = 196.09092
$ws.Range("L6").Value = 600
$ws.Range("M6").Value = -83.09092000000001
$ws.Range("N6").Value = -826

$ws.Range("H68").Value = 830.5454999999999
$ws.Range("I68").Value = 610.1818
$ws.Range("J68").Value = 1050.909
$ws.Range("K68").Value = 1830.5454
$ws.Range("L68").Value = 3152.727
$ws.Range("M68").Value = -1019.5454
$ws.Range("N68").Value = -4774.727000000001

$ws.Range("H71").Value = 830.5454999999999
$ws.Range("I71").Value = 610.1818
$ws.Range("J71").Value = 1050.909
$ws.Range("K71").Value = 5491.6362
$ws.Range("L71").Value = 9458.181
$ws.Range("M71").Value = -1435.6362
$ws.Range("N71").Value = -17570.181

$ws.Range("H131").Value = 12823464
$ws.Range("I131").Value = 55015
$ws.Range("J131").Value = 13159475
$ws.Range("K131").Value = 165045
$ws.Range("L131").Value = 39478425
$ws.Range("M131").Value = -160005
$ws.Range("N131").Value = -39488505

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 929192.25
$ws.Range("I3").Value = 1202569.9
$ws.Range("J3").Value = 17933.334
$ws.Range("K3").Value = 1202569.9
$ws.Range("L3").Value = 17933.334
$ws.Range("M3").Value = -1202453.9
$ws.Range("N3").Value = -18165.334

$ws.Range("H32").Value = 27000
$ws.Range("J32").Value = 27000
$ws.Range("L32").Value = 27000
$ws.Range("N32").Value = -27592

$ws.Range("H45").Value = 38992.332
$ws.Range("J45").Value = 38992.332
$ws.Range("L45").Value = 38992.332
$ws.Range("N45").Value = -40110.332

$ws.Range("H51").Value = 36499.7
$ws.Range("J51").Value = 36499.7
$ws.Range("L51").Value = 36499.7
$ws.Range("N51").Value = -37517.7

$ws.Range("H109").Value = 28285
$ws.Range("J109").Value = 28285
$ws.Range("L109").Value = 28285
$ws.Range("N109").Value = -30365

$ws.Range("H132").Value = 2229.7856
$ws.Range("I132").Value = 1477.2222
$ws.Range("J132").Value = 3584.4
$ws.Range("K132").Value = 4431.6666
$ws.Range("L132").Value = 10753.2
$ws.Range("M132").Value = -1901.6666
$ws.Range("N132").Value = -15813.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 5918.4614
$ws.Range("I9").Value = 454.2
$ws.Range("K9").Value = 454.2
$ws.Range("M9").Value = -230.2

$ws.Range("H81").Value = 33511
$ws.Range("J81").Value = 33511
$ws.Range("L81").Value = 33511
$ws.Range("N81").Value = -35507

$ws.Range("H84").Value = 33511
$ws.Range("J84").Value = 33511
$ws.Range("L84").Value = 100533
$ws.Range("N84").Value = -110517

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 18000
$ws.Range("J11").Value = 10000
$ws.Range("L11").Value = 10000
$ws.Range("N11").Value = -10284

$ws.Range("H123").Value = 22085.584
$ws.Range("J123").Value = 22085.584
$ws.Range("L123").Value = 22085.584
$ws.Range("N123").Value = -31885.584

$ws.Range("H132").Value = 933.9434
$ws.Range("I132").Value = 702.1951
$ws.Range("J132").Value = 1725.75
$ws.Range("K132").Value = 2106.5853
$ws.Range("L132").Value = 5177.25
$ws.Range("M132").Value = 423.4146999999998
$ws.Range("N132").Value = -10237.25

$ws.Range("H140").Value = 50214.5
$ws.Range("J140").Value = 50214.5
$ws.Range("L140").Value = 50214.5
$ws.Range("N140").Value = -60574.5
